# Swap the data rows (row 2 <-> row 3) on the "Quellen Langform 2" sheet,
# and make that sheet the active/selected tab with a new active cell (A7).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Quellen Langform")
$ws2 = $wb.Worksheets.Item("Quellen Langform 2")

# --- Capture current values + row heights before overwriting anything ---
$a2 = $ws2.Cells.Item(2, 1).Value2
$b2 = $ws2.Cells.Item(2, 2).Value2
$a3 = $ws2.Cells.Item(3, 1).Value2
$b3 = $ws2.Cells.Item(3, 2).Value2

$rowHeight2 = $ws2.Rows.Item(2).RowHeight
$rowHeight3 = $ws2.Rows.Item(3).RowHeight

# --- Reset formatting on the four cells so the re-applied wrap-text
#     format is rebuilt cleanly for each destination, instead of the old
#     per-position formatting leaking through ---
$ws2.Cells.Item(2, 1).Clear() | Out-Null
$ws2.Cells.Item(2, 2).Clear() | Out-Null
$ws2.Cells.Item(3, 1).Clear() | Out-Null
$ws2.Cells.Item(3, 2).Clear() | Out-Null

# --- Write the swapped values: old row 3 -> row 2, old row 2 -> row 3 ---
$ws2.Cells.Item(2, 1).Value2 = $a3
$ws2.Cells.Item(2, 2).Value2 = $b3
$ws2.Cells.Item(3, 1).Value2 = $a2
$ws2.Cells.Item(3, 2).Value2 = $b2

# --- Re-apply the wrap-text formatting that belongs with the data now in
#     each row (row 2 / "Luetke" row had no border-styled column A, row 3
#     / "Hertel" row has wrap text on both A and B) ---
$ws2.Cells.Item(3, 1).WrapText = $true
$ws2.Cells.Item(2, 2).WrapText = $true
$ws2.Cells.Item(3, 2).WrapText = $true

# --- Swap the row heights to match the data that now lives there ---
$ws2.Rows.Item(2).RowHeight = $rowHeight3
$ws2.Rows.Item(3).RowHeight = $rowHeight2

# --- "Quellen Langform 2" becomes the active sheet/tab (was "Quellen
#     Langform" before); the first sheet loses tabSelected ---
$ws2.Activate()
$ws2.Select()

# Update the active selection: sheet 1 keeps its old selection, sheet 2's
# selection moves to A7
$ws2.Range("A7").Select()
